$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.22299422634783
$ws.Range("C2").Value = 8.996801150313807
$ws.Range("D2").Value = 11.51655450499791
$ws.Range("F2").Value = 27.68518895197605
$ws.Range("G2").Value = 24.4950714765633
$ws.Range("H2").Value = 13.2709238882142
$ws.Range("J2").Value = 11.18997596994233
$ws.Range("M2").Value = 16.98254418702167
$ws.Range("O2").Value = 19.62102681839587
# Row 3
$ws.Range("B3").Value = 12.56644524088066
$ws.Range("C3").Value = 8.556860847598086
$ws.Range("D3").Value = 11.49442813186268
$ws.Range("F3").Value = 27.78898571170069
$ws.Range("G3").Value = 24.61942599621873
$ws.Range("H3").Value = 13.33361956417977
$ws.Range("J3").Value = 11.22772105426296
$ws.Range("M3").Value = 16.74929109801742
$ws.Range("O3").Value = 19.72805949970147
# Row 4
$ws.Range("B4").Value = 12.14505131518723
$ws.Range("C4").Value = 8.274091327190952
$ws.Range("D4").Value = 11.4826438371689
$ws.Range("F4").Value = 27.86123444975495
$ws.Range("G4").Value = 24.70724173409427
$ws.Range("H4").Value = 13.3748039584227
$ws.Range("J4").Value = 11.25266641127358
$ws.Range("M4").Value = 16.60609102043317
$ws.Range("O4").Value = 19.79936914200351
# Row 5
$ws.Range("B5").Value = 11.96891064808684
$ws.Range("C5").Value = 8.155787229578483
$ws.Range("D5").Value = 11.47829810556364
$ws.Range("F5").Value = 27.89280912608766
$ws.Range("G5").Value = 24.74588626205982
$ws.Range("H5").Value = 13.39226277770262
$ws.Range("J5").Value = 11.26327686367837
$ws.Range("M5").Value = 16.54780021121669
$ws.Range("O5").Value = 19.82982981826488
# Row 6
$ws.Range("B6").Value = 11.93940142864352
$ws.Range("C6").Value = 8.135960783161632
$ws.Range("D6").Value = 11.47760416511131
$ws.Range("F6").Value = 27.89818061827813
$ws.Range("G6").Value = 24.75247513063336
$ws.Range("H6").Value = 13.39520261160692
$ws.Range("J6").Value = 11.26506559941904
$ws.Range("M6").Value = 16.53812670738845
$ws.Range("O6").Value = 19.83497230056387
# Row 7
$ws.Range("B7").Value = 12.14269346034029
$ws.Range("C7").Value = 8.27250812603396
$ws.Range("D7").Value = 11.48258337653303
$ws.Range("F7").Value = 27.8616516537243
$ws.Range("G7").Value = 24.70775136351034
$ws.Range("H7").Value = 13.37503667811965
$ws.Range("J7").Value = 11.25280770567992
$ws.Range("M7").Value = 16.60530455033774
$ws.Range("O7").Value = 19.79977427703234
# Row 8
$ws.Range("B8").Value = 13.0005040128598
$ws.Range("C8").Value = 8.847792317918007
$ws.Range("D8").Value = 11.50855364066426
$ws.Range("F8").Value = 27.71920508489159
$ws.Range("G8").Value = 24.53555641195859
$ws.Range("H8").Value = 13.29198305788585
$ws.Range("J8").Value = 11.20262317302723
$ws.Range("M8").Value = 16.90215049830577
$ws.Range("O8").Value = 19.65676852913016
# Row 9
$ws.Range("B9").Value = 14.53142921042059
$ws.Range("C9").Value = 9.871788862693197
$ws.Range("D9").Value = 11.57360837497604
$ws.Range("F9").Value = 27.50782920905308
$ws.Range("G9").Value = 24.28981268114628
$ws.Range("H9").Value = 13.15046575369858
$ws.Range("J9").Value = 11.11825309926145
$ws.Range("M9").Value = 17.48180768227678
$ws.Range("O9").Value = 19.420899695221
# Row 10
$ws.Range("B10").Value = 15.55727383450511
$ws.Range("C10").Value = 10.55656199752476
$ws.Range("D10").Value = 11.62977551268819
$ws.Range("F10").Value = 27.39445548481372
$ws.Range("G10").Value = 24.16657864059884
$ws.Range("H10").Value = 13.05952905197012
$ws.Range("J10").Value = 11.06482605139433
$ws.Range("M10").Value = 17.90275967959853
$ws.Range("O10").Value = 19.27505048627597
# Row 11
$ws.Range("B11").Value = 16.00145195809728
$ws.Range("C11").Value = 10.85279524589559
$ws.Range("D11").Value = 11.65708746123939
$ws.Range("F11").Value = 27.35207709060716
$ws.Range("G11").Value = 24.12320615840502
$ws.Range("H11").Value = 13.02099513156283
$ws.Range("J11").Value = 11.04237937907692
$ws.Range("M11").Value = 18.09249128794276
$ws.Range("O11").Value = 19.21471820992363
# Row 12
$ws.Range("B12").Value = 16.16635257413016
$ws.Range("C12").Value = 10.96273562453577
$ws.Range("D12").Value = 11.667677653849
$ws.Range("F12").Value = 27.33735872578984
$ws.Range("G12").Value = 24.10862440724357
$ws.Range("H12").Value = 13.00681134133299
$ws.Range("J12").Value = 11.03414657459353
$ws.Range("M12").Value = 18.16402793168515
$ws.Range("O12").Value = 19.19274188596158
# Row 13
$ws.Range("B13").Value = 16.1309859724966
$ws.Range("C13").Value = 10.93915803576303
$ws.Range("D13").Value = 11.66538593846292
$ws.Range("F13").Value = 27.34046935738214
$ws.Range("G13").Value = 24.11168263426146
$ws.Range("H13").Value = 13.0098479095519
$ws.Range("J13").Value = 11.03590776734086
$ws.Range("M13").Value = 18.14863594015565
$ws.Range("O13").Value = 19.19743609536148
# Row 14
$ws.Range("B14").Value = 16.01508490503962
$ws.Range("C14").Value = 10.86188513770668
$ws.Range("D14").Value = 11.65795379007298
$ws.Range("F14").Value = 27.35083952028831
$ws.Range("G14").Value = 24.12196948566812
$ws.Range("H14").Value = 13.01982003821035
$ws.Range("J14").Value = 11.04169670408258
$ws.Range("M14").Value = 18.09838319652217
$ws.Range("O14").Value = 19.21289273463464
# Row 15
$ws.Range("B15").Value = 15.94366063454744
$ws.Range("C15").Value = 10.81426090517711
$ws.Range("D15").Value = 11.65343348046026
$ws.Range("F15").Value = 27.3573648719856
$ws.Range("G15").Value = 24.12851093030762
$ws.Range("H15").Value = 13.02598143171018
$ws.Range("J15").Value = 11.04527740500859
$ws.Range("M15").Value = 18.06755982489649
$ws.Range("O15").Value = 19.22247385257715
# Row 16
$ws.Range("B16").Value = 15.5277873008951
$ws.Range("C16").Value = 10.53689147673029
$ws.Range("D16").Value = 11.62802558854858
$ws.Range("F16").Value = 27.39741067477653
$ws.Range("G16").Value = 24.16967001521667
$ws.Range("H16").Value = 13.06210441318724
$ws.Range("J16").Value = 11.06633038566496
$ws.Range("M16").Value = 17.89031988337855
$ws.Range("O16").Value = 19.27911483080932
# Row 17
$ws.Range("B17").Value = 15.26685200384555
$ws.Range("C17").Value = 10.36279103966306
$ws.Range("D17").Value = 11.61288593483285
$ws.Range("F17").Value = 27.42433798297424
$ws.Range("G17").Value = 24.19818265681264
$ws.Range("H17").Value = 13.08499094170091
$ws.Range("J17").Value = 11.07972157733968
$ws.Range("M17").Value = 17.78109636120737
$ws.Range("O17").Value = 19.31540660269703
# Row 18
$ws.Range("B18").Value = 15.1146557216384
$ws.Range("C18").Value = 10.26121731288458
$ws.Range("D18").Value = 11.6043441310913
$ws.Range("F18").Value = 27.4406910350488
$ws.Range("G18").Value = 24.2157755827276
$ws.Range("H18").Value = 13.09842134000462
$ws.Range("J18").Value = 11.08759866978861
$ws.Range("M18").Value = 17.7181113816662
$ws.Range("O18").Value = 19.33684654381811
# Row 19
$ws.Range("B19").Value = 15.06276398745256
$ws.Range("C19").Value = 10.22658096631616
$ws.Range("D19").Value = 11.60148071838446
$ws.Range("F19").Value = 27.4463762671616
$ws.Range("G19").Value = 24.2219366034574
$ws.Range("H19").Value = 13.10301442263474
$ws.Range("J19").Value = 11.09029574126082
$ws.Range("M19").Value = 17.69675961825741
$ws.Range("O19").Value = 19.34420276980049
# Row 20
$ws.Range("B20").Value = 15.29484824594806
$ws.Range("C20").Value = 10.38147324237729
$ws.Range("D20").Value = 11.61448042255652
$ws.Range("F20").Value = 27.42138192938238
$ws.Range("G20").Value = 24.19502380888892
$ws.Range("H20").Value = 13.08252702711617
$ws.Range("J20").Value = 11.07827796648112
$ws.Range("M20").Value = 17.79274063084405
$ws.Range("O20").Value = 19.31148467797907
# Row 21
$ws.Range("B21").Value = 16.04921788794035
$ws.Range("C21").Value = 10.88464304794921
$ws.Range("D21").Value = 11.66013011370134
$ws.Range("F21").Value = 27.34775741907664
$ws.Range("G21").Value = 24.11889784837476
$ws.Range("H21").Value = 13.01687989902899
$ws.Range("J21").Value = 11.03998909838935
$ws.Range("M21").Value = 18.11315252160697
$ws.Range("O21").Value = 19.20832908658277
# Row 22
$ws.Range("B22").Value = 16.52298482325946
$ws.Range("C22").Value = 11.20044284670741
$ws.Range("D22").Value = 11.69140614745988
$ws.Range("F22").Value = 27.30739120300811
$ws.Range("G22").Value = 24.07989162348059
$ws.Range("H22").Value = 12.97635516963984
$ws.Range("J22").Value = 11.01652299534735
$ws.Range("M22").Value = 18.3207246925468
$ws.Range("O22").Value = 19.14598553833352
# Row 23
$ws.Range("B23").Value = 16.27190737322582
$ws.Range("C23").Value = 11.03310005997697
$ws.Range("D23").Value = 11.67458355362223
$ws.Range("F23").Value = 27.32822396912406
$ws.Range("G23").Value = 24.0997211435802
$ws.Range("H23").Value = 12.99776602267795
$ws.Range("J23").Value = 11.02890469957853
$ws.Range("M23").Value = 18.21012586645024
$ws.Range("O23").Value = 19.17879342055651
# Row 24
$ws.Range("B24").Value = 15.28219793335998
$ws.Range("C24").Value = 10.37303162861478
$ws.Range("D24").Value = 11.61375904905868
$ws.Range("F24").Value = 27.42271564494756
$ws.Range("G24").Value = 24.19644818554006
$ws.Range("H24").Value = 13.08364011371818
$ws.Range("J24").Value = 11.0789300675095
$ws.Range("M24").Value = 17.78747684683999
$ws.Range("O24").Value = 19.31325598844014
# Row 25
$ws.Range("B25").Value = 14.13424092805421
$ws.Range("C25").Value = 9.606389527948746
$ws.Range("D25").Value = 11.55452053447345
$ws.Range("F25").Value = 27.55768331318398
$ws.Range("G25").Value = 24.34631501451069
$ws.Range("H25").Value = 13.15046575369858
$ws.Range("J25").Value = 11.11825309926145
$ws.Range("M25").Value = 17.32561601625176
$ws.Range("O25").Value = 19.420899695221

Write-Host "Applied 216 cell updates"